# Scheduled-runner style refresh of the Ultros_Profits sheets.
# Updates the market-derived columns (currentAveragePrice / NQ / HQ,
# LevePriceNQ / HQ, LeveProfitNQ / HQ -> columns H..N) for the specific
# leve rows whose source market data changed. Values only; no formulas,
# no structural changes.

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H28").Value = 1283.2354
$ws.Range("I28").Value = 1262.1428
$ws.Range("K28").Value = 1262.1428
$ws.Range("M28").Value = -777.1428000000001

$ws.Range("H53").Value = 211.1923
$ws.Range("I53").Value = 86
$ws.Range("J53").Value = 277.47058
$ws.Range("K53").Value = 86
$ws.Range("L53").Value = 277.47058
$ws.Range("M53").Value = 551
$ws.Range("N53").Value = -1551.47058

$ws.Range("H137").Value = 3853.8518
$ws.Range("I137").Value = 3070.4546
$ws.Range("J137").Value = 4392.4375
$ws.Range("K137").Value = 9211.363799999999
$ws.Range("L137").Value = 13177.3125
$ws.Range("M137").Value = -6661.363799999999
$ws.Range("N137").Value = -18277.3125

$ws.Range("H138").Value = 3019.8333
$ws.Range("I138").Value = 1910.6666
$ws.Range("J138").Value = 3322.3333
$ws.Range("K138").Value = 5731.9998
$ws.Range("L138").Value = 9966.999899999999
$ws.Range("M138").Value = -591.9997999999996
$ws.Range("N138").Value = -20246.9999

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H74").Value = 3164.6086
$ws.Range("I74").Value = 3132.7646
$ws.Range("K74").Value = 3132.7646
$ws.Range("M74").Value = -2258.7646

$ws.Range("H77").Value = 3164.6086
$ws.Range("I77").Value = 3132.7646
$ws.Range("K77").Value = 15663.823
$ws.Range("M77").Value = -11295.823

$ws.Range("H97").Value = 1438.6316
$ws.Range("I97").Value = 869.2692
$ws.Range("J97").Value = 2672.25
$ws.Range("K97").Value = 869.2692
$ws.Range("L97").Value = 2672.25
$ws.Range("M97").Value = -373.2692
$ws.Range("N97").Value = -3664.25

$ws.Range("H122").Value = 4704.8623
$ws.Range("J122").Value = 5818.091
$ws.Range("L122").Value = 17454.273
$ws.Range("N122").Value = -22354.273

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

# --- BSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H86").Value = 25051794
$ws.Range("I86").Value = 62626576
$ws.Range("J86").Value = 1936.9166
$ws.Range("K86").Value = 62626576
$ws.Range("L86").Value = 1936.9166
$ws.Range("M86").Value = -62625453
$ws.Range("N86").Value = -4182.9166

$ws.Range("H89").Value = 25051794
$ws.Range("I89").Value = 62626576
$ws.Range("J89").Value = 1936.9166
$ws.Range("K89").Value = 313132880
$ws.Range("L89").Value = 9684.583000000001
$ws.Range("M89").Value = -313127264
$ws.Range("N89").Value = -20916.583

$ws.Range("H99").Value = 46604.55
$ws.Range("J99").Value = 85137.164
$ws.Range("L99").Value = 85137.164
$ws.Range("N99").Value = -88133.164

$ws.Range("H105").Value = 5700
$ws.Range("I105").Value = 5000
$ws.Range("J105").Value = 5875
$ws.Range("K105").Value = 5000
$ws.Range("L105").Value = 5875
$ws.Range("M105").Value = -3253
$ws.Range("N105").Value = -9369

$ws.Range("H134").Value = 2744.0908
$ws.Range("I134").Value = 2256
$ws.Range("K134").Value = 6768
$ws.Range("M134").Value = -4233

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 8800
$ws.Range("J16").Value = 8800
$ws.Range("L16").Value = 8800
$ws.Range("N16").Value = -9374

$ws.Range("H31").Value = 4120.9287
$ws.Range("I31").Value = 3175.25
$ws.Range("J31").Value = 4499.2
$ws.Range("K31").Value = 3175.25
$ws.Range("L31").Value = 4499.2
$ws.Range("M31").Value = -2880.25
$ws.Range("N31").Value = -5089.2

$ws.Range("H34").Value = 4120.9287
$ws.Range("I34").Value = 3175.25
$ws.Range("J34").Value = 4499.2
$ws.Range("K34").Value = 3175.25
$ws.Range("L34").Value = 4499.2
$ws.Range("M34").Value = -2973.25
$ws.Range("N34").Value = -4903.2

$ws.Range("H58").Value = 2770.1
$ws.Range("I58").Value = 1874.5555
$ws.Range("K58").Value = 1874.5555
$ws.Range("M58").Value = -1671.5555

$ws.Range("H107").Value = 9916.272000000001
$ws.Range("I107").Value = 340.8
$ws.Range("K107").Value = 340.8
$ws.Range("M107").Value = 1579.2

$ws.Range("H113").Value = 8800
$ws.Range("J113").Value = 8800
$ws.Range("L113").Value = 8800
$ws.Range("N113").Value = -13140

$ws.Range("H136").Value = 2770.1
$ws.Range("I136").Value = 1874.5555
$ws.Range("K136").Value = 5623.666499999999
$ws.Range("M136").Value = -3073.666499999999

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 1152.2142
$ws.Range("I5").Value = 1052.5834
$ws.Range("J5").Value = 1750
$ws.Range("K5").Value = 3157.7502
$ws.Range("L5").Value = 5250
$ws.Range("M5").Value = -3045.7502
$ws.Range("N5").Value = -5474

$ws.Range("H135").Value = 1152.2142
$ws.Range("I135").Value = 1052.5834
$ws.Range("J135").Value = 1750
$ws.Range("K135").Value = 9473.250599999999
$ws.Range("L135").Value = 15750
$ws.Range("M135").Value = -6938.250599999999
$ws.Range("N135").Value = -20820

# --- LTW -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 2800.5
$ws.Range("I7").Value = 2760.6
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2760.6
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2648.6
$ws.Range("N7").Value = -3224

$ws.Range("H61").Value = 1825.25
$ws.Range("I61").Value = 1825.25
$ws.Range("K61").Value = 1825.25
$ws.Range("M61").Value = -1623.25

$ws.Range("H93").Value = 1086.6666
$ws.Range("I93").Value = 904
$ws.Range("K93").Value = 904
$ws.Range("M93").Value = 344

$ws.Range("H113").Value = 1825.25
$ws.Range("I113").Value = 1825.25
$ws.Range("K113").Value = 1825.25
$ws.Range("M113").Value = 344.75

$ws.Range("H126").Value = 2800.5
$ws.Range("I126").Value = 2760.6
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8281.799999999999
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5811.799999999999
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 2405.4194
$ws.Range("I132").Value = 2296.6538
$ws.Range("K132").Value = 6889.9614
$ws.Range("M132").Value = -4359.9614

$ws.Range("H136").Value = 2238.5908
$ws.Range("I136").Value = 2199.2104
$ws.Range("J136").Value = 2488
$ws.Range("K136").Value = 6597.6312
$ws.Range("L136").Value = 7464
$ws.Range("M136").Value = -4047.6312
$ws.Range("N136").Value = -12564

# --- WVR -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H114").Value = 39999.5
$ws.Range("J114").Value = 39999.5
$ws.Range("L114").Value = 39999.5
$ws.Range("N114").Value = -48677.5

$ws.Range("H122").Value = 1795.5714
$ws.Range("I122").Value = 1626
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 4878
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2428
$ws.Range("N122").Value = -16900
